# "1st changes of mifos to finflux"
#
# On the "Repayment schedule" sheet, insert a new blank column immediately
# before the existing "Late" column (column N). The previous N/O/P columns
# (Late / heading / Outstanding) shift right to O/P/Q. The new column is
# given the same width as its left neighbour (the "In Advance" column).
#
# Also switch the active sheet/tab from "Transactions" to "Repayment
# schedule", with cell S6 selected there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at position N (14); the old N..P columns
# (and their data/styles) shift right to O..Q.
$ws.Columns.Item(14).Insert() | Out-Null

# Give the freshly inserted column the same width as column M (to its left).
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select S6 - this also
# clears the previously active selection/tab on the "Transactions" sheet.
$ws.Activate() | Out-Null
$ws.Range("S6").Select() | Out-Null
